$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source scrape re-ran and re-ordered several fixtures that share the
# same kickoff date/time, plus picked up one new match. For each pair of
# rows below, the match details (home/away teams, scores, odds, odds
# timestamps and match URL -- columns F:V) were swapped between the two
# rows, while the row's own index/metadata (columns A:E) stayed put.

$swapPairs = @(
    @(84, 85),
    @(96, 97),
    @(111, 112),
    @(132, 133)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("F$r1`:V$r1")
    $range2 = $ws.Range("F$r2`:V$r2")
    $tmp = $range1.Value2
    $range1.Value = $range2.Value2
    $range2.Value = $tmp
}

# A new fixture (Goztepe vs Eyupspor) was appended as row 143. Copy the
# formatting of the last existing row down first so the new row picks up
# the same cell styles (bordered/bold index column, date-formatted column
# E), then fill in its actual values.
$ws.Range("A142:V142").Copy($ws.Range("A143:V143"))

$ws.Range("A143").Value = 142
$ws.Range("B143").Value = "turkey"
$ws.Range("C143").Value = "1-lig"
$ws.Range("D143").Value = "2023-2024"
$ws.Range("E143").Value = 45280.75
$ws.Range("F143").Value = "Goztepe"
$ws.Range("G143").Value = 5
$ws.Range("H143").Value = "Eyupspor"
$ws.Range("I143").Value = 1
$ws.Range("J143").Value = 3.41
$ws.Range("K143").Value = "10/12/2023 17:12"
$ws.Range("L143").Value = 3.59
$ws.Range("M143").Value = "20/12/2023 17:57"
$ws.Range("N143").Value = 2.89
$ws.Range("O143").Value = "10/12/2023 17:12"
$ws.Range("P143").Value = 3.21
$ws.Range("Q143").Value = "20/12/2023 17:57"
$ws.Range("R143").Value = 2.31
$ws.Range("S143").Value = "10/12/2023 17:12"
$ws.Range("T143").Value = 2.2
$ws.Range("U143").Value = "20/12/2023 17:57"
$ws.Range("V143").Value = "https://www.betexplorer.com/football/turkey/1-lig/goztepe-eyupspor/Iwtl7KHO/"
